$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Historical data revisions (rows 212-230): corrections to
# "Nb nouvelles admissions a l'hopital" (D) and "Patients COVID-19
# hospitalises hors SI" (G) ---
$ws.Range("D212").Value = 2
$ws.Range("G212").Value = 9

$ws.Range("G213").Value = 10

$ws.Range("G214").Value = 9

$ws.Range("G215").Value = 9

$ws.Range("D216").Value = 1
$ws.Range("G216").Value = 9

$ws.Range("G217").Value = 11

$ws.Range("G218").Value = 12

$ws.Range("G219").Value = 14

$ws.Range("G220").Value = 14

$ws.Range("G221").Value = 15

$ws.Range("G222").Value = 18

$ws.Range("D223").Value = 1
$ws.Range("G223").Value = 18

$ws.Range("G224").Value = 15

$ws.Range("G225").Value = 16

$ws.Range("G226").Value = 17

$ws.Range("G227").Value = 21

$ws.Range("G228").Value = 24

$ws.Range("G229").Value = 27

$ws.Range("G230").Value = 30

# --- Rows 231-237: corrections to "Nb nouveaux cas positifs" (C),
# admissions (D), SI patients (E), hospitalises hors SI (G) and
# nouvelles sorties (I) ---
$ws.Range("C231").Value = 223
$ws.Range("D231").Value = 6
$ws.Range("G231").Value = 36

$ws.Range("C232").Value = 302
$ws.Range("D232").Value = 13
$ws.Range("G232").Value = 49

$ws.Range("C233").Value = 307
$ws.Range("G233").Value = 51

$ws.Range("C234").Value = 395
$ws.Range("D234").Value = 13
$ws.Range("G234").Value = 63

$ws.Range("D235").Value = 19
$ws.Range("G235").Value = 82

$ws.Range("C236").Value = 193
$ws.Range("D236").Value = 6
$ws.Range("G236").Value = 85

$ws.Range("C237").Value = 307
$ws.Range("D237").Value = 2
$ws.Range("E237").Value = 6
$ws.Range("G237").Value = 82
$ws.Range("I237").Value = 3

# --- Row 238: fill in new day's data (2020-09-21) ---
$ws.Range("C238").Value = 19
$ws.Range("D238").Value = 0
$ws.Range("E238").Value = 8
$ws.Range("F238").Value = 1
$ws.Range("G238").Value = 80
$ws.Range("I238").Value = 0
$ws.Range("L238").Value = "0"
$ws.Range("M238").Value = "0"

# --- Update the view selection to reflect where the user ended up ---
$ws.Range("B234").Select()
$excel.ActiveWindow.ScrollRow = 234
$ws.Range("O238").Select()
